$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column E header: date 2024-03-03 (serial 45354), formatted as a date
# (set format first so Excel doesn't mint a throwaway custom numFmt before
# the explicit one is applied)
$ws.Range("E1").NumberFormat = "mm-dd-yy"
$ws.Range("E1").Value = "3/3/2024"

# New column E data values
$ws.Range("E2").Value = 1314
$ws.Range("E3").Value = 1355
$ws.Range("E4").Value = 1330
$ws.Range("E5").Value = 3054
$ws.Range("E6").Value = 155

# Size column E to fit its new contents (best-fit column width)
$ws.Columns("E").ColumnWidth = 9.5

# Update selection to match where the user ended up after entering the data
$ws.Range("E6").Select() | Out-Null
